$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 with the new server data
# Column order per header row: A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
# Set values in the same order the new shared strings were added (ServerID, then IP, then ID/Name)
$ws.Range("B2").Value = "000106001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "MasterServer_1"
$ws.Range("C2").Value = "MasterServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# A2/B2 already carry the text style; make sure C2 matches it too
$ws.Range("C2").NumberFormat = "@"

# Update the active selection to reflect H3 being the active cell
$ws.Range("H3").Select()
